$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 2444, 20.40344619750977, 3.812446355819702, 14, 49),
    @(2, 2444, 29.22617149353027, 12.87995433807373, 13, 93),
    @(3, 2444, 29.05969619750977, 15.14703750610352, 13, 74),
    @(4, 2444, 19.48077011108398, 2.365187406539917, 15, 27),
    @(1, 2444, 44.13541793823242, 32.65696716308594, 15, 140),
    @(2, 2444, 135.1153869628906, 68.97154235839844, 23, 254)
)

$startRow = 50
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
